# Apply updated hourly generation/capacity regression output values
# (Coef./Std.Err./t/P>|t|/CI bounds/coef_pos) per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("H2").Value = 0.1131786086946073

# Row 3
$ws.Range("B3").Value = 0.1218613636993578
$ws.Range("H3").Value = 0.235039972393965

# Row 4
$ws.Range("B4").Value = 0.1375932818981117
$ws.Range("H4").Value = 0.2507718905927189

# Row 5
$ws.Range("B5").Value = 0.04878831074351646
$ws.Range("H5").Value = 0.1619669194381237

# Row 6
$ws.Range("B6").Value = 0.0343649643497862
$ws.Range("C6").Value = 0.002158421334054636
$ws.Range("D6").Value = 5.202212969087995
$ws.Range("E6").Value = 0.02335440618728778
$ws.Range("F6").Value = 0.03013220212614183
$ws.Range("G6").Value = 0.03859772657343079
$ws.Range("H6").Value = 0.1475435730443935

# Row 7
$ws.Range("B7").Value = 0.02359404533206398
$ws.Range("C7").Value = 0.002039481093661821
$ws.Range("D7").Value = 4.688935704704987
$ws.Range("E7").Value = 0.005642390157575884
$ws.Range("F7").Value = 0.01959212637235279
$ws.Range("G7").Value = 0.02759596429177455
$ws.Range("H7").Value = 0.1367726540266712

# Row 8
$ws.Range("B8").Value = 0.01941021750111112
$ws.Range("C8").Value = 0.001461175360400451
$ws.Range("D8").Value = 1.748050355832855
$ws.Range("E8").Value = 0.003696734588197441
$ws.Range("F8").Value = 0.01654507921841213
$ws.Range("G8").Value = 0.02227535578380962
$ws.Range("H8").Value = 0.1325888261957184

# Row 9
$ws.Range("B9").Value = 0.0166185622620552
$ws.Range("C9").Value = 0.001610544705972097
$ws.Range("D9").Value = 1.564815158547401
$ws.Range("E9").Value = 0.006230055483430611
$ws.Range("F9").Value = 0.01346043771353193
$ws.Range("G9").Value = 0.01977668681057825
$ws.Range("H9").Value = 0.1297971709566625

# Row 10
$ws.Range("B10").Value = 0.01714509198296621
$ws.Range("C10").Value = 0.001629357016901622
$ws.Range("D10").Value = 1.51464351573503
$ws.Range("E10").Value = 0.006553274137775946
$ws.Range("F10").Value = 0.01395018109562983
$ws.Range("G10").Value = 0.02034000287030291
$ws.Range("H10").Value = 0.1303237006775735

# Row 11
$ws.Range("B11").Value = 0.02984110674673522
$ws.Range("H11").Value = 0.1430197154413425

# Row 12
$ws.Range("B12").Value = 0.05641442708176776
$ws.Range("H12").Value = 0.169593035776375

# Row 13
$ws.Range("B13").Value = 0.07580404521182262
$ws.Range("H13").Value = 0.1889826539064299

# Row 14
$ws.Range("B14").Value = 0.08310892516889606
$ws.Range("H14").Value = 0.1962875338635033

# Row 15
$ws.Range("B15").Value = 0.0903514308826778
$ws.Range("H15").Value = 0.203530039577285

# Row 16
$ws.Range("B16").Value = 0.09514804186755348
$ws.Range("H16").Value = 0.2083266505621607

# Row 17
$ws.Range("B17").Value = 0.09795295822564463
$ws.Range("H17").Value = 0.2111315669202519

# Row 18
$ws.Range("B18").Value = -0.1131786086946073

# Row 19
$ws.Range("B19").Value = 0.100736314709637
$ws.Range("H19").Value = 0.2139149234042443

# Row 20
$ws.Range("B20").Value = 0.1025697862965416
$ws.Range("H20").Value = 0.2157483949911489

# Row 21
$ws.Range("B21").Value = 0.1063911053297102
$ws.Range("H21").Value = 0.2195697140243175

# Row 22
$ws.Range("B22").Value = 0.1124193032281466
$ws.Range("H22").Value = 0.2255979119227538

# Row 23
$ws.Range("B23").Value = 0.1128603529604227
$ws.Range("H23").Value = 0.22603896165503

# Row 24
$ws.Range("B24").Value = 0.1200310554705017
$ws.Range("C24").Value = 0.006643996483768277
$ws.Range("D24").Value = 1566163970759.371
$ws.Range("E24").Value = 0.02919966064840934
$ws.Range("F24").Value = 0.1069690528182446
$ws.Range("G24").Value = 0.1330930581227591
$ws.Range("H24").Value = 0.2332096641651089

# Row 25
$ws.Range("B25").Value = 0.1211818174852429
$ws.Range("C25").Value = 0.006686943135642252
$ws.Range("D25").Value = 1509254108728.727
$ws.Range("E25").Value = 0.03270129652224789
$ws.Range("F25").Value = 0.1080363451759315
$ws.Range("G25").Value = 0.134327289794554
$ws.Range("H25").Value = 0.2343604261798502

# Row 26
$ws.Range("B26").Value = 0.1248624877034587
$ws.Range("C26").Value = 0.006532891869858073
$ws.Range("D26").Value = -126401847077.5133
$ws.Range("E26").Value = 0.03369208752984119
$ws.Range("F26").Value = 0.1120239103561863
$ws.Range("G26").Value = 0.1377010650507302
$ws.Range("H26").Value = 0.2380410963980659

# Row 27
$ws.Range("B27").Value = 0.1283076976407428
$ws.Range("C27").Value = 0.006600492221340341
$ws.Range("D27").Value = 29.99144521546262
$ws.Range("E27").Value = 0.04451991747656985
$ws.Range("F27").Value = 0.1153362128591081
$ws.Range("G27").Value = 0.1412791824223772
$ws.Range("H27").Value = 0.2414863063353501

# Row 28
$ws.Range("B28").Value = 0.1321520855296657
$ws.Range("C28").Value = 0.006628135629196968
$ws.Range("D28").Value = 2031777523345.456
$ws.Range("E28").Value = 0.07209140588759202
$ws.Range("F28").Value = 0.1191340085878469
$ws.Range("G28").Value = 0.145170162471484
$ws.Range("H28").Value = 0.245330694224273

# Row 29
$ws.Range("B29").Value = 0.0206156902622914
$ws.Range("C29").Value = 0.001483712743870398
$ws.Range("D29").Value = 2.424025198025625
$ws.Range("E29").Value = 0.01290611937003616
$ws.Range("F29").Value = 0.01770573049253632
$ws.Range("G29").Value = 0.02352565003204622
$ws.Range("H29").Value = 0.1337942989568986
